# Apply the target edits to the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: "3.28×104" -> "3.28×10^4"
$ws.Range("C2").Value = "3.28×10^4"

# F2: "1.00×102" -> "1.00×10^2"
$ws.Range("F2").Value = "1.00×10^2"

# C4: "110.80" -> "110.8"
$ws.Range("C4").Value = "110.8"

# E4: "TU/L" -> "U/L"
$ws.Range("E4").Value = "U/L"

# F5: "<0,500" -> "<0.500"
$ws.Range("F5").Value = "<0.500"
